$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Logs sheet: append new row 7 (new mail log entry) ---
$logs.Range("A7").Value = "Bel klant"
$logs.Range("B7").Value = "mailmind.test@zohomail.eu"
$logs.Range("D7").Value = "Klantenservice / Contact"
$logs.Range("F7").Value = "2025-08-30 19:24:30"
$logs.Range("G7").Value = "Nee"
$logs.Range("H7").Value = "Ja"
$logs.Range("I7").Value = "Nee"
$logs.Range("J7").Value = "Nee"

# --- Dashboard sheet: append new category row 5 ---
$dash.Range("A5").Value = "Klantenservice / Contact"
$dash.Range("B5").Value = 1

# --- Extend conditional formatting ranges on Logs from row 6 to row 7 ---
$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $oldRange = $logs.Range($col + "2:" + $col + "6")
    $newRange = $logs.Range($col + "2:" + $col + "7")
    $fc = $oldRange.FormatConditions
    for ($i = 1; $i -le $fc.Count; $i++) {
        $fc.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Extend the Dashboard bar chart's category/value series to include row 5 ---
$chartObj = $dash.ChartObjects(1)
$chart = $chartObj.Chart
$ser = $chart.SeriesCollection(1)
$ser.XValues = "='Dashboard'!`$A`$2:`$A`$5"
$ser.Values = "='Dashboard'!`$B`$2:`$B`$5"
